# "Add files via upload" — small content tweaks to the team's plan sheet:
#   - extend the final summary note
#   - fill in the "completion status" column for the last few rows
#   - scroll position / window size tweaks left by the author's session

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "完成情况" (completion status) cells.
$ws.Range("C104").Value = "只完成了一部分"
$ws.Range("C105").Value = "完成"
$ws.Range("C106").Value = "完成"
$ws.Range("C107").Value = "完成"
$ws.Range("C108").Value = "只完成了一部分"
$ws.Range("C109").Value = "完成"

# Expand the closing summary remark.
$ws.Range("A110").Value = "总结：做后台的要加把劲"

# Cosmetic: restore the scroll position / window size from the author's
# last editing session (best effort — view-only state).
$win = $excel.ActiveWindow
$win.ScrollRow = 95
$win.ScrollColumn = 1

try {
    $win.Width = 19095
    $win.Height = 12210
} catch {
    # window sizing may be a no-op in headless hosts; ignore
}
